$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 13.0107
$ws.Range("E6").Value = 12.3096
$ws.Range("C7").Value = -12.7127
$ws.Range("A8").Value = -21.24430000000001
$ws.Range("E9").Value = 10.30659999999999
$ws.Range("A10").Value = -20.46139999999997
$ws.Range("E10").Value = 11.7257
$ws.Range("A12").Value = -22.68260000000003
$ws.Range("B13").Value = 5.954799999999997
$ws.Range("A18").Value = -22.60780000000003
$ws.Range("C20").Value = -14.8628
